# Update the "registered family" (row 4) and "subsistence-allowance
# recipient family" (row 5) counts for 2015-2021 (columns E-K) on the
# "გორი" sheet with refreshed figures from the ministry data source.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: რეგისტრირებული ოჯახი (registered families), years 2015-2021
$ws.Range("E4").Value = 20280
$ws.Range("F4").Value = 14521
$ws.Range("G4").Value = 13573
$ws.Range("H4").Value = 12935
$ws.Range("I4").Value = 12260
$ws.Range("J4").Value = 13048
$ws.Range("K4").Value = 13789

# Row 5: საარსებო შემწეობის მიმღები ოჯახი (subsistence-allowance
# recipient families), years 2015-2021
$ws.Range("E5").Value = 6623
$ws.Range("F5").Value = 6797
$ws.Range("G5").Value = 5842
$ws.Range("H5").Value = 5541
$ws.Range("I5").Value = 4719
$ws.Range("J5").Value = 5766
$ws.Range("K5").Value = 6805
